$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2"  = -0.608236435814144
    "C2"  = -1.81872946919787
    "B3"  = 0.202516267794457
    "C3"  = 1.20638425134075
    "B4"  = 1.29447715192983
    "C4"  = 4.43505114580974
    "B5"  = -0.275869033621741
    "C5"  = -0.71456314590315
    "B6"  = -0.176039326449044
    "C6"  = 0.921891106884248
    "B7"  = -0.946062823975874
    "C7"  = -1.49122835747614
    "B8"  = -2.18966187662582
    "C8"  = 3.28370102827252
    "B9"  = -0.17430788600259
    "C9"  = 2.53819688998371
    "B10" = -0.29356529269873
    "C10" = 3.57240875105334
    "B11" = 0.878130621548514
    "C11" = 0.228346786173341
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
